$d = $word.ActiveDocument

# 1. Title: "Alex Alvarez Gárciga" -> "Alex Alvarez Garciga"
$d.Content.Find.Execute("Alvarez Gárciga", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Alvarez Garciga", 2)

# 2. Address line: "Homestead, FL, USA." -> "Miami, FL, USA."
$d.Content.Find.Execute("Homestead, FL, USA.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Miami, FL, USA.", 2)

# 3. "In love with React.js" -> "In love with React"
$d.Content.Find.Execute("In love with React.js", $true, $false, $false, $false, $false,
                         $true, 1, $false, "In love with React", 2)

# 4. "Universidad de las Ciencias Informáticas" -> "Universidad de las Ciencias Informaticas"
$d.Content.Find.Execute("Universidad de las Ciencias Informáticas", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Universidad de las Ciencias Informaticas", 2)
